# Applies the score-correction edit described in the commit
# "tournament closure exception handeled" to the match scoresheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Temba Bavuma(C) / Tamim Iqbal
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 'LBW'
$ws.Range("E2").Value = ' Mahedi Hasan'
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 9
$ws.Range("N2").Value = ' Kagiso Rabada'

# Row 3: Quinton de Kock / Liton Das
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Shoriful Islam'
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 3
$ws.Range("N3").Value = ' Anrich Nortje'

# Row 4: Rassie Va der Dussen / Shakib Al Hasan
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 'LBW'
$ws.Range("K4").Value = 13
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 'NOT OUT'
$ws.Range("N4").Value = ' '

# Row 5: Aiden Markram / Mushfiqur Rahim
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 'Caught'
$ws.Range("N5").Value = ' Dwaine Pretorius'

# Row 6: David Miller / Mahmudulla(C)
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 'LBW'
$ws.Range("E6").Value = ' Taskin Ahmed'
$ws.Range("K6").Value = 6
$ws.Range("M6").Value = 'Caught'
$ws.Range("N6").Value = ' Dwaine Pretorius'

# Row 7: Reeza Hendricks / Afif Hossain
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = ' Shamim Hossain'
$ws.Range("K7").Value = 11
$ws.Range("L7").Value = 8
$ws.Range("M7").Value = '* NOT OUT'
$ws.Range("N7").Value = ' '

# Row 8: Dwaine Pretorius / Shamim Hossain
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 'Bowled'
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = ' '
$ws.Range("N8").Value = ' '

# Row 9: Kagiso Rabada / Mahedi Hasan
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 'Bowled'
$ws.Range("E9").Value = ' Shamim Hossain'
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = ' '
$ws.Range("N9").Value = ' '

# Row 10: Keshav Maharaj / Taskin Ahmed
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 'NOT OUT'
$ws.Range("E10").Value = ' '
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ' '

# Row 11: Anrich Nortje / Mustafizur Rahman
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 'Bowled'
$ws.Range("E11").Value = ' Shoriful Islam'
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ' '
$ws.Range("N11").Value = ' '

# Row 12: Tabraiz Shamsi / Shoriful Islam
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 'LBW'
$ws.Range("E12").Value = ' Shoriful Islam'
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ' '
$ws.Range("N12").Value = ' '

# Row 16: Innings totals
$ws.Range("A16").Value = 58
$ws.Range("C16").Value = '''6.0'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = 36
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = '''5.5'
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = 35

# Row 21: Bowling figures (left: Mustafizur Rahman, right: Anrich Nortje)
$ws.Range("A21").Value = 'Mustafizur Rahman'
$ws.Range("B21").Value = '''1.0'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 15
$ws.Range("K21").Value = '''1.0'
$ws.Range("K21").Style = "Normal"
$ws.Range("L21").Value = 9
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 9

# Row 22: Bowling figures (left: Taskin Ahmed, right: Keshav Maharaj)
$ws.Range("A22").Value = 'Taskin Ahmed'
$ws.Range("B22").Value = '''1.0'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 9
$ws.Range("K22").Value = '''1.0'
$ws.Range("K22").Style = "Normal"
$ws.Range("L22").Value = 9
$ws.Range("N22").Value = 9

# Row 23: Bowling figures (left: Mahedi Hasan, right: Kagiso Rabada)
$ws.Range("A23").Value = 'Mahedi Hasan'
$ws.Range("B23").Value = '''1.0'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 16
$ws.Range("E23").Value = 16
$ws.Range("K23").Value = '''1.0'
$ws.Range("K23").Style = "Normal"
$ws.Range("L23").Value = 6
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 6

# Row 24: Bowling figures (left: Shamim Hossain, right: Dwaine Pretorius)
$ws.Range("A24").Value = 'Shamim Hossain'
$ws.Range("B24").Value = '''1.0'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 7
$ws.Range("K24").Value = '''1.0'
$ws.Range("K24").Style = "Normal"
$ws.Range("L24").Value = 11
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 11

# Row 25: Bowling figures (left: Shoriful Islam, right: Tabraiz Shamsi)
$ws.Range("A25").Value = 'Shoriful Islam'
$ws.Range("B25").Value = '''2.0'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 5.5
$ws.Range("K25").Value = '''1.5'
$ws.Range("K25").Style = "Normal"
$ws.Range("L25").Value = 27
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 18

Write-Output "Applied all scoresheet corrections"
